# Commit: "switched back to realistic annual vehicle mileages"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "average annual mileage per vehicle (U.S.) [mi]" column (P)
# back to more realistic per-income-group figures.
$ws.Range("P2").Value = 13250
$ws.Range("P3").Value = 14250
$ws.Range("P4").Value = 13750
# P5 (13000) is unchanged
$ws.Range("P6").Value = 12500
$ws.Range("P7").Value = 11500

# Move the active selection from P11 back to B3 (and, implicitly, scroll
# the view back so column A is visible again instead of being scrolled to H1).
$ws.Activate()
$ws.Range("B3").Select()
